$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '45.348.59'
$ws.Range("E2").Value = '  -0.45%  '
# Row 3
$ws.Range("D3").Value = '2.367.85'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.58%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.88'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.27%  '
# Row 8
$ws.Range("E8").Value = '  +0.02%  '
# Row 9
$ws.Range("E9").Value = '  -3.31%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.73'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.82%  '
# Row 11
$ws.Range("E11").Value = '  -1.44%  '
# Row 12
$ws.Range("E12").Value = '  -2.24%  '
# Row 13
$ws.Range("E13").Value = '  +1.08%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.979'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.91%  '
# Row 15
$ws.Range("D15").Value = '2.727.95'
$ws.Range("E15").Value = '  -0.43%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.30'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.76%  '
# Row 17
$ws.Range("D17").Value = '2.364.49'
$ws.Range("E17").Value = '  -0.78%  '
# Row 18
$ws.Range("D18").Value = '45.375.79'
$ws.Range("E18").Value = '  +0.16%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.26'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +9.48%  '
# Row 20
$ws.Range("E20").Value = '  -1.81%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.20'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.56%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.14'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.53%  '
# Row 23
$ws.Range("E23").Value = '  -0.40%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.83'
$ws.Range("D24").ClearFormats()
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.46%  '
# Row 26
$ws.Range("E26").Value = '  +0.02%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.32%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.26'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.89%  '
# Row 29
$ws.Range("E29").Value = '  -1.63%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0971'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.24%  '
# Row 31
$ws.Range("E31").Value = '  -2.80%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.78'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -5.08%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '166.37'
$ws.Range("D33").ClearFormats()
# Row 34
$ws.Range("E34").Value = '  -3.17%  '
# Row 35
$ws.Range("E35").Value = '  -2.12%  '
# Row 36
$ws.Range("E36").Value = '  +0.52%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.70'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.49%  '
# Row 38
$ws.Range("E38").Value = '  +9.06%  '
# Row 39
$ws.Range("E39").Value = '  +0.33%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.07%  '
# Row 41
$ws.Range("E41").Value = '  -3.57%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.17'
$ws.Range("D42").ClearFormats()
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '70.03'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.09%  '
# Row 44
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.226'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.90%  '
# Row 45
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.02%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.76'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -7.63%  '
# Row 47
$ws.Range("D47").Value = '1.808.67'
$ws.Range("E47").Value = '  +9.30%  '
# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.83'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.72%  '
# Row 49
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.13'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.93%  '
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.59'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.18%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.19'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.03%  '
